$d = $word.ActiveDocument

# 1. "construct expectations about how unfamiliar talkers will sound." -> "constructs expectations about how that talker will produce speech in the future."
$d.Content.Find.Execute(
    "construct expectations about how unfamiliar talkers will sound.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "constructs expectations about how that talker will produce speech in the future.",
    2)

# 2. "However, this process still requires cognitive resources. In this experiment, we limit the available resources"
#    -> "However, how automatic speech perception adaptation is remains unclear. In this experiment, we limit the available attentional resources"
$d.Content.Find.Execute(
    "However, this process still requires cognitive resources. In this experiment, we limit the available resources",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "However, how automatic speech perception adaptation is remains unclear. In this experiment, we limit the available attentional resources",
    2)

# 3. "suggest a difference between passive and active attention in speech processing. Additionally, the results of this experiment will give insight into how our brains allocate resources"
#    -> "suggest there are limits to the automaticity of speech perception. Additionally, the results of this experiment will provide insight into how our brains allocates attentional resources"
$d.Content.Find.Execute(
    "suggest a difference between passive and active attention in speech processing. Additionally, the results of this experiment will give insight into how our brains allocate resources",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "suggest there are limits to the automaticity of speech perception. Additionally, the results of this experiment will provide insight into how our brains allocates attentional resources",
    2)
